$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows below header and write new table content
$ws.Range("A2:C19").ClearContents()

$ws.Range("A2").Value = 'Jamal Murray'
$ws.Range("B2").Value = 'PG,SG'
$ws.Range("C2").Value = 'Denver Nuggets'
$ws.Range("A3").Value = 'Donovan Mitchell'
$ws.Range("B3").Value = 'PG,SG'
$ws.Range("C3").Value = 'Cleveland Cavaliers'
$ws.Range("A4").Value = 'Malik Beasley'
$ws.Range("B4").Value = 'SG,SF'
$ws.Range("C4").Value = 'Detroit Pistons'
$ws.Range("A5").Value = 'Michael Porter Jr.'
$ws.Range("B5").Value = 'SF,PF'
$ws.Range("C5").Value = 'Denver Nuggets'
$ws.Range("A6").Value = 'Domantas Sabonis'
$ws.Range("B6").Value = 'C'
$ws.Range("C6").Value = 'Sacramento Kings'
$ws.Range("A7").Value = 'Andre Drummond'
$ws.Range("B7").Value = 'C'
$ws.Range("C7").Value = 'Philadelphia 76ers'
$ws.Range("A8").Value = 'Victor Wembanyama'
$ws.Range("B8").Value = 'C'
$ws.Range("C8").Value = 'San Antonio Spurs'
$ws.Range("A9").Value = 'Kristaps Porzingis'
$ws.Range("B9").Value = 'PF,C'
$ws.Range("C9").Value = 'Boston Celtics'
$ws.Range("A10").Value = 'Myles Turner'
$ws.Range("B10").Value = 'C'
$ws.Range("C10").Value = 'Indiana Pacers'
$ws.Range("A11").Value = 'Payton Pritchard'
$ws.Range("B11").Value = 'PG'
$ws.Range("C11").Value = 'Boston Celtics'
$ws.Range("A12").Value = 'Josh Hart'
$ws.Range("B12").Value = 'SG,SF,PF'
$ws.Range("C12").Value = 'New York Knicks'
$ws.Range("A13").Value = 'Bradley Beal'
$ws.Range("B13").Value = 'PG,SG,SF'
$ws.Range("C13").Value = 'Phoenix Suns'
$ws.Range("A14").Value = 'Tari Eason'
$ws.Range("B14").Value = 'SF,PF'
$ws.Range("C14").Value = 'Houston Rockets'
$ws.Range("A15").Value = 'De''Andre Hunter'
$ws.Range("B15").Value = 'SF,PF'
$ws.Range("C15").Value = 'Atlanta Hawks'
$ws.Range("A16").Value = 'Dyson Daniels'
$ws.Range("B16").Value = 'PG,SG,SF'
$ws.Range("C16").Value = 'Atlanta Hawks'
$ws.Range("A17").Value = 'Cam Thomas'
$ws.Range("B17").Value = 'SG,SF'
$ws.Range("C17").Value = 'Brooklyn Nets'
$ws.Range("A18").Value = 'Deandre Ayton'
$ws.Range("B18").Value = 'C'
$ws.Range("C18").Value = 'Portland Trail Blazers'

# Remove the now-unused last row (table shrank from 18 to 17 data rows)
$ws.Rows.Item(19).Delete()

